# Add new "MX-DPBX" / "MX-BBX" accessory rows into the panel-accessories
# list on the Italy and Netherlands market sheets (rows are inserted right
# before the trailing "Wg" / "Accessories" marker rows, i.e. above what
# used to be row 12).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Italy", "Netherlands")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Make room for the two new rows above the old row 12 ("Wg") /
    # row 13 ("Accessories"), shifting them down to rows 14/15.
    $ws.Rows.Item(12).Insert()
    $ws.Rows.Item(12).Insert()

    # Copy the formatting (thin border style) of the row above down onto
    # the two freshly inserted rows so they match the rest of the list.
    $ws.Range("A11").Copy()
    $ws.Range("A12:A13").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    $ws.Range("A12").Value = "MX-DPBX"
    $ws.Range("A13").Value = "MX-BBX"

    [void]$ws.Range("A12:A13").Select()
}

# The Netherlands sheet previously had no explicit page setup; give it the
# same portrait orientation used by the other market sheets.
$wsNetherlands = $wb.Worksheets.Item("Netherlands")
$wsNetherlands.PageSetup.Orientation = 1
